$wb = $excel.ActiveWorkbook

# "linear" sheet
$ws1 = $wb.Worksheets.Item("linear")
$ws1.Range("B2").Value = 0.004105376846097051
$ws1.Range("B3").Value = -0.07898817841423227
$ws1.Range("B4").Value = 1.353521011482407

# "non-linear" sheet
$ws2 = $wb.Worksheets.Item("non-linear")
$ws2.Range("B2").Value = 0.02052115074596907
$ws2.Range("B3").Value = 0.01509938526654823
$ws2.Range("B4").Value = 1.379653801527025
$ws2.Range("B5").Value = 0.08066323695490685
$ws2.Range("B6").Value = -0.276795127221192
$ws2.Range("B7").Value = 1.324415593620945
